$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: add a new bulleted paragraph right after the existing bullet
# "Si se deshabilita un rol, ningún usuario podrá acceder con ese rol",
# using the same list style/numbering as its sibling bullets.
# -----------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*deshabilita un rol*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "Los usuarios siguen teniendo ese rol, así en el caso de que se vuelva a habilitar pueden volver a acceder con el mismo."
        break
    }
}

# -----------------------------------------------------------------------
# Change 2: the <w:lastRenderedPageBreak/> marker moves from the run that
# starts the "El precio de la factura..." bullet to the run that holds
# the "Facturas" heading text (it now renders at the top of the new
# page instead of on the first bullet below it).
# -----------------------------------------------------------------------
$facturasPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Facturas") {
        $facturasPara = $p
        break
    }
}

if ($facturasPara -ne $null) {
    $nextPara = $facturasPara.Next()

    # Add the page-break marker right before the "Facturas" text.
    $r1 = $facturasPara.Range
    $xml1 = $r1.WordOpenXML
    if ($xml1 -like "*<w:lastRenderedPageBreak/>*") {
        # already present (nothing to do)
    } else {
        $xml1 = $xml1.Replace("<w:t>Facturas</w:t>", "<w:lastRenderedPageBreak/><w:t>Facturas</w:t>")
        $r1.InsertXML($xml1)
    }

    # Remove the page-break marker from the following bullet paragraph.
    if ($nextPara -ne $null) {
        $r2 = $nextPara.Range
        $xml2 = $r2.WordOpenXML
        if ($xml2 -like "*<w:lastRenderedPageBreak/>*") {
            $xml2 = $xml2.Replace("<w:lastRenderedPageBreak/>", "")
            $r2.InsertXML($xml2)
        }
    }
}
